$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.138.58"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.918.93"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.89"
$ws.Range("E5").Value = "  -3.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5061"
$ws.Range("E7").Value = "  -3.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4041"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08283"
$ws.Range("E9").Value = "  -2.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.114"
$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.30"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.87"
$ws.Range("E12").Value = "  +3.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.422"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.900.30"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.347"
$ws.Range("E15").Value = "  -0.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.69"
$ws.Range("E17").Value = "  -2.82%  "

$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06485"
$ws.Range("E19").Value = "  -3.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.61"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.977"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.232.25"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.195"
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.39"
$ws.Range("E26").Value = "  +5.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.116.99"
$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.96"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.394"
$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.04"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.130"
$ws.Range("E31").Value = "  +2.45%  "

$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.993"
$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.804"
$ws.Range("E34").Value = "  +5.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02460"
$ws.Range("E35").Value = "  -1.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.403"
$ws.Range("E36").Value = "  +4.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06472"
$ws.Range("E37").Value = "  -1.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2163"
$ws.Range("E38").Value = "  -2.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.746"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.190"
$ws.Range("E40").Value = "  -3.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6420"
$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.40"
$ws.Range("E42").Value = "  -4.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.220"
$ws.Range("E43").Value = "  -1.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9982"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.183"
$ws.Range("E45").Value = "  +4.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.26"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6005"
$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.646"
$ws.Range("E48").Value = "  -2.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.89"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.219"
$ws.Range("E50").Value = "  -2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.11"
$ws.Range("E51").Value = "  -0.85%  "
